$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bets")
$ws.Rows(25).Copy()
$ws.Rows(26).Insert(-4121)
$ws.Rows(26).Copy()
$ws.Rows(27).Insert(-4121)

$ws.Range("A26").Value = 25
$ws.Range("B26").Value = 45220
$ws.Range("E26").Value = 930
$ws.Range("I26").Value = "JDG"
$ws.Range("J26").Value = "GANA SERIE"

$ws.Range("A27").Value = 26
$ws.Range("B27").Value = 45220
$ws.Range("E27").Value = 4
$ws.Range("I27").Value = "GEN"
$ws.Range("J27").Value = "GANA SERIE"
Write-Output "done"
